# Updates the crypto price table (cols D=Price, E=Volume(1h)) to the latest
# scrape, and swaps the Cronos/EnergySwap rows (48/49) to match new ranking order.
# D-column values that look like plain numbers are forced back to Text via
# NumberFormat "@" before the write, mirroring how the source data keeps prices
# (e.g. "16.50", "0.000009615") as literal strings instead of numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.247.84'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.827.39'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.27'
$ws.Range("E5").Value = '  -2.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5984'
$ws.Range("E6").Value = '  -4.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06955'
$ws.Range("E8").Value = '  -5.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2752'
$ws.Range("E9").Value = '  -4.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.23'
$ws.Range("E10").Value = '  -6.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07611'
$ws.Range("D12").Value = '1.832.27'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.739'
$ws.Range("E13").Value = '  -4.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6247'
$ws.Range("E14").Value = '  -6.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009615'
$ws.Range("E15").Value = '  -7.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '78.41'
$ws.Range("E16").Value = '  -3.76%  '
$ws.Range("D17").Value = '28.737.66'
$ws.Range("E17").Value = '  -1.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.596'
$ws.Range("E18").Value = '  -10.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.06'
$ws.Range("E19").Value = '  -6.67%  '
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.53'
$ws.Range("E21").Value = '  -5.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.842'
$ws.Range("E22").Value = '  -6.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.006'
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '156.28'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.939'
$ws.Range("E25").Value = '  -6.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1283'
$ws.Range("E26").Value = '  -3.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.50'
$ws.Range("E27").Value = '  -4.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.441'
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06322'
$ws.Range("E29").Value = '  -11.33%  '
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.828'
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.747'
$ws.Range("E32").Value = '  -7.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.719'
$ws.Range("E33").Value = '  -5.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.086'
$ws.Range("E34").Value = '  -5.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6453'
$ws.Range("E35").Value = '  -8.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.542'
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.749'
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01750'
$ws.Range("E38").Value = '  -4.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.569'
$ws.Range("E39").Value = '  -3.44%  '
$ws.Range("D40").Value = '1.151.08'
$ws.Range("E40").Value = '  -6.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8891'
$ws.Range("E41").Value = '  -5.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.005'
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").Value = '1.984.51'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.44'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.90'
$ws.Range("E45").Value = '  -4.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000115'
$ws.Range("E46").Value = '  -2.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.595'
$ws.Range("E47").Value = '  -5.37%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05533'
$ws.Range("E48").Value = '  -1.95%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.408'
$ws.Range("E49").Value = '  -5.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4553'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.406'
$ws.Range("E51").Value = '  -7.66%  '
